$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.085.96"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.118.47"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.66"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5198"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4461"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.08"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09360"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.183"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.39"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.587"
$ws.Range("E13").Value = "  +6.30%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.128.51"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.963"
$ws.Range("E15").Value = "  +3.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.68"
$ws.Range("E16").Value = "  +3.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001165"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.57"
$ws.Range("E19").Value = "  +4.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06695"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.290"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.118.37"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.74"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.14"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.534"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.52"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "134.15"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.791"
$ws.Range("E31").Value = "  +10.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1057"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.267"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.615"
$ws.Range("E34").Value = "  +7.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.964"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.78"
$ws.Range("E36").Value = "  +6.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02626"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06869"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.7092"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.69"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2248"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.326"
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6859"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.62"
$ws.Range("E44").Value = "  +2.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.385"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.005"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.265"
$ws.Range("E47").Value = "  +7.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.634"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000349"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.227"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "83.35"
$ws.Range("E51").Value = "  +1.52%  "
